$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5071.4287
$ws.Range("I62").Value = 3600
$ws.Range("K62").Value = 3600
$ws.Range("M62").Value = -2976
$ws.Range("H64").Value = 2950
$ws.Range("H65").Value = 5071.4287
$ws.Range("I65").Value = 3600
$ws.Range("K65").Value = 18000
$ws.Range("M65").Value = -14880
$ws.Range("H67").Value = 2950
$ws.Range("H116").Value = 463032.38
$ws.Range("I116").Value = 1001940.5
$ws.Range("J116").Value = 13942.25
$ws.Range("K116").Value = 1001940.5
$ws.Range("L116").Value = 13942.25
$ws.Range("M116").Value = -998498.5
$ws.Range("N116").Value = -20826.25
$ws.Range("H123").Value = 43280
$ws.Range("J123").Value = 43280
$ws.Range("L123").Value = 43280
$ws.Range("N123").Value = -53080
$ws.Range("H132").Value = 107811.52
$ws.Range("I132").Value = 127096.625
$ws.Range("K132").Value = 381289.875
$ws.Range("M132").Value = -378759.875

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1291.25
$ws.Range("I2").Value = 1291.25
$ws.Range("K2").Value = 1291.25
$ws.Range("M2").Value = -1178.25
$ws.Range("H45").Value = 2366.923
$ws.Range("I45").Value = 1270
$ws.Range("K45").Value = 1270
$ws.Range("M45").Value = -893
$ws.Range("H63").Value = 12594491
$ws.Range("I63").Value = 15392533
$ws.Range("J63").Value = 3300
$ws.Range("K63").Value = 15392533
$ws.Range("L63").Value = 3300
$ws.Range("M63").Value = -15391847
$ws.Range("N63").Value = -4672
$ws.Range("H66").Value = 12594491
$ws.Range("I66").Value = 15392533
$ws.Range("J66").Value = 3300
$ws.Range("K66").Value = 76962665
$ws.Range("L66").Value = 16500
$ws.Range("M66").Value = -76959233
$ws.Range("N66").Value = -23364
$ws.Range("H88").Value = 8336595.5
$ws.Range("I88").Value = 22224688
$ws.Range("J88").Value = 3740
$ws.Range("K88").Value = 22224688
$ws.Range("L88").Value = 3740
$ws.Range("M88").Value = -22224282
$ws.Range("N88").Value = -4552
$ws.Range("H91").Value = 8336595.5
$ws.Range("I91").Value = 22224688
$ws.Range("J91").Value = 3740
$ws.Range("K91").Value = 22224688
$ws.Range("L91").Value = 3740
$ws.Range("M91").Value = -22223284
$ws.Range("N91").Value = -6548
$ws.Range("H116").Value = 1291.25
$ws.Range("I116").Value = 1291.25
$ws.Range("K116").Value = 1291.25
$ws.Range("M116").Value = 1002.75

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1291.25
$ws.Range("I3").Value = 1291.25
$ws.Range("K3").Value = 1291.25
$ws.Range("M3").Value = -1177.25
$ws.Range("H64").Value = 424.75
$ws.Range("J64").Value = 499.66666
$ws.Range("L64").Value = 499.66666
$ws.Range("N64").Value = -949.66666
$ws.Range("H67").Value = 424.75
$ws.Range("J67").Value = 499.66666
$ws.Range("L67").Value = 499.66666
$ws.Range("N67").Value = -2059.66666
$ws.Range("H105").Value = 2441.3333
$ws.Range("I105").Value = 2393.3333
$ws.Range("J105").Value = 2633.3333
$ws.Range("K105").Value = 2393.3333
$ws.Range("L105").Value = 2633.3333
$ws.Range("M105").Value = -646.3332999999998
$ws.Range("N105").Value = -6127.3333

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3379.1035
$ws.Range("I31").Value = 1135.1111
$ws.Range("J31").Value = 7051.091
$ws.Range("K31").Value = 1135.1111
$ws.Range("L31").Value = 7051.091
$ws.Range("M31").Value = -840.1111000000001
$ws.Range("N31").Value = -7641.091
$ws.Range("H34").Value = 3379.1035
$ws.Range("I34").Value = 1135.1111
$ws.Range("J34").Value = 7051.091
$ws.Range("K34").Value = 1135.1111
$ws.Range("L34").Value = 7051.091
$ws.Range("M34").Value = -933.1111000000001
$ws.Range("N34").Value = -7455.091
$ws.Range("H99").Value = 4230.5713
$ws.Range("I99").Value = 2513.7778
$ws.Range("J99").Value = 7320.8
$ws.Range("K99").Value = 2513.7778
$ws.Range("L99").Value = 7320.8
$ws.Range("M99").Value = -1015.7778
$ws.Range("N99").Value = -10316.8
$ws.Range("H126").Value = 4230.5713
$ws.Range("I126").Value = 2513.7778
$ws.Range("J126").Value = 7320.8
$ws.Range("K126").Value = 7541.3334
$ws.Range("L126").Value = 21962.4
$ws.Range("M126").Value = -5071.3334
$ws.Range("N126").Value = -26902.4
$ws.Range("H139").Value = 47660
$ws.Range("J139").Value = 47660
$ws.Range("L139").Value = 47660
$ws.Range("N139").Value = -57940

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1378.8125
$ws.Range("I5").Value = 626.3333
$ws.Range("J5").Value = 3034.2666
$ws.Range("K5").Value = 1878.9999
$ws.Range("L5").Value = 9102.799800000001
$ws.Range("M5").Value = -1766.9999
$ws.Range("N5").Value = -9326.799800000001
$ws.Range("H98").Value = 197.8
$ws.Range("J98").Value = 277
$ws.Range("L98").Value = 831
$ws.Range("N98").Value = -3827
$ws.Range("H122").Value = 2751.0845
$ws.Range("J122").Value = 3691.8
$ws.Range("L122").Value = 33226.2
$ws.Range("N122").Value = -38126.2
$ws.Range("H131").Value = 11112020
$ws.Range("I131").Value = 27778690
$ws.Range("J131").Value = 907.03705
$ws.Range("K131").Value = 83336070
$ws.Range("L131").Value = 2721.11115
$ws.Range("M131").Value = -83331030
$ws.Range("N131").Value = -12801.11115
$ws.Range("H135").Value = 1378.8125
$ws.Range("I135").Value = 626.3333
$ws.Range("J135").Value = 3034.2666
$ws.Range("K135").Value = 5636.9997
$ws.Range("L135").Value = 27308.3994
$ws.Range("M135").Value = -3101.9997
$ws.Range("N135").Value = -32378.3994
$ws.Range("H141").Value = 7534.647
$ws.Range("I141").Value = 6776.5557
$ws.Range("J141").Value = 8387.5
$ws.Range("K141").Value = 20329.6671
$ws.Range("L141").Value = 25162.5
$ws.Range("M141").Value = -15149.6671
$ws.Range("N141").Value = -35522.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7129.2915
$ws.Range("I70").Value = 6442.9414
$ws.Range("J70").Value = 8796.143
$ws.Range("K70").Value = 6442.9414
$ws.Range("L70").Value = 8796.143
$ws.Range("M70").Value = -6172.9414
$ws.Range("N70").Value = -9336.143
$ws.Range("H73").Value = 7129.2915
$ws.Range("I73").Value = 6442.9414
$ws.Range("J73").Value = 8796.143
$ws.Range("K73").Value = 6442.9414
$ws.Range("L73").Value = 8796.143
$ws.Range("M73").Value = -5506.9414
$ws.Range("N73").Value = -10668.143
$ws.Range("H80").Value = 20836080
$ws.Range("I80").Value = 62501750
$ws.Range("J80").Value = 3244.5
$ws.Range("K80").Value = 62501750
$ws.Range("L80").Value = 3244.5
$ws.Range("M80").Value = -62500752
$ws.Range("N80").Value = -5240.5
$ws.Range("H83").Value = 20836080
$ws.Range("I83").Value = 62501750
$ws.Range("J83").Value = 3244.5
$ws.Range("K83").Value = 312508750
$ws.Range("L83").Value = 16222.5
$ws.Range("M83").Value = -312503758
$ws.Range("N83").Value = -26206.5
$ws.Range("H116").Value = 29000
$ws.Range("J116").Value = 29000
$ws.Range("L116").Value = 29000
$ws.Range("N116").Value = -38178
$ws.Range("H132").Value = 2216.0488
$ws.Range("I132").Value = 1110.8928
$ws.Range("J132").Value = 4596.385
$ws.Range("K132").Value = 3332.6784
$ws.Range("L132").Value = 13789.155
$ws.Range("M132").Value = -802.6784000000002
$ws.Range("N132").Value = -18849.155
$ws.Range("H136").Value = 26989.941
$ws.Range("J136").Value = 26989.941
$ws.Range("L136").Value = 80969.823
$ws.Range("N136").Value = -86069.823

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5802.5
$ws.Range("I7").Value = 4000
$ws.Range("J7").Value = 7154.375
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 7154.375
$ws.Range("M7").Value = -3888
$ws.Range("N7").Value = -7378.375
$ws.Range("H126").Value = 5802.5
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 7154.375
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 21463.125
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -26403.125
$ws.Range("H128").Value = 42721.5
$ws.Range("J128").Value = 42721.5
$ws.Range("L128").Value = 42721.5
$ws.Range("N128").Value = -52681.5
$ws.Range("H136").Value = 3021.6667
$ws.Range("I136").Value = 1756.375
$ws.Range("J136").Value = 4033.9
$ws.Range("K136").Value = 5269.125
$ws.Range("L136").Value = 12101.7
$ws.Range("M136").Value = -2719.125
$ws.Range("N136").Value = -17201.7
